$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.564.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.457.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.836.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.482.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.772"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.553.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.967.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.695.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  -1.08%  "
